# The commit "unify the conception of DataNode, DataTable, Entity" renames
# the single worksheet that used to be called "Property1" so that it now
# reflects the unified "DataNode" concept. The rest of the original XML
# diff (fileVersion/build numbers, absPath, window geometry, xr/xr2/xr3
# revision GUIDs, style-locale names, default row height, etc.) is just
# incidental metadata written by whichever Excel build/platform the author
# happened to resave the file with - there is no cell data or structural
# change behind it. The one visible, intentional UI change captured in the
# diff besides the rename is the selected range on the sheet, which we
# reproduce as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet: "Property1" -> "DataNode"
$ws.Name = "DataNode"

# Restore the frozen-pane view and update the current selection to match
# the author's saved selection (A9:N35) while keeping the same active cell.
$ws.Range("A9:N35").Select()
